$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data content of row 4 and row 5 (all fields that differ
# between the two records), including moving the "På tall" public comment
# from row 4 to row 5.

# --- Row 4 gets row 5's original values ---
$ws.Range("A4").Value = 111470245
$ws.Range("B4").Value = 96348
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("Q4").Value = 554481.1995954363
$ws.Range("R4").Value = 7003291.317192273
$ws.Range("Z4").Value = "14:41"
$ws.Range("AB4").Value = "14:41"
$ws.Range("AC4").ClearContents()

# --- Row 5 gets row 4's original values ---
$ws.Range("A5").Value = 111471797
$ws.Range("B5").Value = 77515
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 554597.2688619854
$ws.Range("R5").Value = 7003280.616068945
$ws.Range("Z5").Value = "15:49"
$ws.Range("AB5").Value = "15:49"
$ws.Range("AC5").Value = "På tall"
